$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing checkout time for Wednesday (row 15, column D) - 5:00 PM
$ws.Range("D15").Value = 0.708333333333333

# Re-apply the print area, which (as in the real Excel/LO interop flow) adds
# a duplicate _xlnm.Print_Area_0 defined name alongside the existing one
$ws.Names.Add("_xlnm.Print_Area_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# Move the active selection to E5
[void]$ws.Range("E5").Select()
